# Capstone Hour Tracker - add the 10/14/2023 (row 23) time entry, with
# its accomplishment / next-step notes, and select the next log cell
# (H24) ready for the following entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New time-log row (row 23) ---------------------------------------
$ws.Range("B23").Value = 45213                      # Date: 10/14/2023
$ws.Range("C23").Value = 0.53125                    # Start Time: 12:45 PM
$ws.Range("C23").NumberFormat = "[$-409]h:mm\ AM/PM;@"
$ws.Range("D23").Value = 0.82291666666666663        # End Time: 7:45 PM
$ws.Range("D23").NumberFormat = "h:mm AM/PM"

# Session notes
$ws.Range("G23").Value = "Got the file explorer opening files! I also got it so files generated by npm i are shown to the user. Also you can select a specific file using tabs in the code editor"
$ws.Range("H23").Value = "Next up is the terminal. Also I need to deal with technical debt soon if not now"

# Row grew tall to fit the wrapped notes, same as the other logged days
$ws.Rows.Item(23).RowHeight = 60

# --- Move the active selection down to the next empty log row --------
[void]$ws.Range("H24").Select()
